$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$B2Text = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f986427c580>),
                ('model',
                 LogisticRegression(class_weight='balanced', l1_ratio=0.3,
                                    max_iter=1000, penalty='elasticnet',
                                    random_state=42, solver='saga'))])
'@
$ws.Range("B2").Value = $B2Text

$C2Text = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f98e1369400>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.3, 'model__class_weight': 'balanced'}
'@
$ws.Range("C2").Value = $C2Text

$B3Text = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f98e1369ac0>),
                ('model',
                 LogisticRegression(class_weight='balanced', l1_ratio=0.01,
                                    max_iter=1000, penalty='elasticnet',
                                    random_state=42, solver='saga'))])
'@
$ws.Range("B3").Value = $B3Text

$C3Text = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f98b0e6ca60>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.01, 'model__class_weight': 'balanced'}
'@
$ws.Range("C3").Value = $C3Text

$B4Text = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f98e0ec3c10>),
                ('model',
                 LogisticRegression(class_weight='balanced', l1_ratio=0.7,
                                    max_iter=1000, penalty='elasticnet',
                                    random_state=42, solver='saga'))])
'@
$ws.Range("B4").Value = $B4Text

$C4Text = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f98643259d0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.7, 'model__class_weight': 'balanced'}
'@
$ws.Range("C4").Value = $C4Text

$B5Text = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f98643252b0>),
                ('model',
                 LogisticRegression(class_weight='balanced', l1_ratio=0.7,
                                    max_iter=1000, penalty='elasticnet',
                                    random_state=42, solver='saga'))])
'@
$ws.Range("B5").Value = $B5Text

$C5Text = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f98e0cef4f0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.7, 'model__class_weight': 'balanced'}
'@
$ws.Range("C5").Value = $C5Text

$B6Text = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 LogisticRegression(class_weight='balanced', l1_ratio=0.1,
                                    max_iter=1000, penalty='elasticnet',
                                    random_state=42, solver='saga'))])
'@
$ws.Range("B6").Value = $B6Text

$C6Text = @'
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.1, 'model__class_weight': 'balanced'}
'@
$ws.Range("C6").Value = $C6Text

$ws.Range("D6").Value = 0.548528176628863

$ws.Range("G6").Value = 0.6999969365560763

$ws.Range("H6").Value = 0.6269841269841269

$J6Text = @'
[0 0 1 1 1 0 0 0 1 1 1 1 0 1 1 0 0 0 0 1 0 1 1 1]
'@
$ws.Range("J6").Value = $J6Text
